# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on row 2 of the
# zh-cn and de-de worksheets to reflect the new report generation times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 15:18:24"
$wsZhCn.Range("H2").Value = "2016-03-23 15:18:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 15:18:28"
$wsDeDe.Range("H2").Value = "2016-03-23 15:18:57"
